{"js": "// Office.js (Word JavaScript API) script.\n// Replaces the English table-header labels with their Chinese counterparts,\n// per the target diff:\n//   \"Logistic Regression Model for Cars\" -> \"\u8f66\u4fe9\u6570\u636e\u903b\u8f91\u56de\u5f52\u6a21\u578b\"\n//   \"Odds Ratio\"                          -> \"\u4f18\u52bf\u6bd4\"\n//   \"Std. error\"                          -> \"\u6807\u51c6\u8bef\u5dee\"\n//   \"z\"                                   -> \"Z\u5206\u6570\"\n//   \"p-value\"                             -> \"P\u503c\"\n//   \"95% CI\"                               -> \"95\"\n\nconst replacements = [\n  [\"Logistic Regression Model for Cars\", \"\u8f66\u4fe9\u6570\u636e\u903b\u8f91\u56de\u5f52\u6a21\u578b\", false],\n  [\"Odds Ratio\", \"\u4f18\u52bf\u6bd4\", false],\n  [\"Std. error\", \"\u6807\u51c6\u8bef\u5dee\", false],\n  [\"z\", \"Z\u5206\u6570\", true],\n  [\"p-value\", \"P\u503c\", false],\n  [\"95% CI\", \"95\", false],\n];\n\nfor (const [find, replace, wholeWord] of replacements) {\n  const results = context.document.body.search(find, {\n    matchCase: true,\n    matchWholeWord: wholeWord,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${find}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @{\n    \"Logistic Regression Model for Cars\" = \"\u8f66\u4fe9\u6570\u636e\u903b\u8f91\u56de\u5f52\u6a21\u578b\"\n    \"Odds Ratio\" = \"\u4f18\u52bf\u6bd4\"\n    \"Std. error\" = \"\u6807\u51c6\u8bef\u5dee\"\n    \"z\" = \"Z\u5206\u6570\"\n    \"p-value\" = \"P\u503c\"\n    \"95% CI\" = \"95\"\n}\n\nforeach ($find in $pairs.Keys) {\n    $replace = $pairs[$find]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($find, $false, $true, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
